$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows 2-9 each hold one species observation record. This edit
# re-orders those records: the record that used to sit on row 3
# (TaxonId 5260, "Lateritticka" / Postia lateritia) and the record that
# used to sit on row 9 (TaxonId 498, "Liten sotlav" / Acolium karelicum)
# are moved up to directly follow the header row, and the remaining
# records (originally on rows 2,4,5,6,7,8) shift down to fill rows 4-9,
# keeping their original relative order.
#
# Rather than physically moving whole rows (which would disturb the
# sparse cell layout of unrelated, unchanged "blank" cells such as the
# empty-string placeholders in columns I/AT/AY or the boolean
# AD/AE/AG flags), only the cells whose value actually differs between
# the old and new layout are written here, cell by cell.

# New value for each touched cell, keyed by "<col><row>".
# A value of $null means the cell must be cleared entirely (used for the
# "Antal substrat" / "Substrat-beskrivning" cells AN/AO, which exist only
# on two of the eight rows and move together with their owning record).
$updates = [ordered]@{
    "A2"  = 1561978
    "B2"  = 89610
    "E2"  = 5260
    "F2"  = "Lateritticka"
    "G2"  = "Postia lateritia"
    "H2"  = "Renvall"
    "Q2"  = 783440.5104734434
    "R2"  = 7532071.274571465

    "A3"  = 175486
    "B3"  = 76862
    "E3"  = 498
    "F3"  = "Liten sotlav"
    "G3"  = "Acolium karelicum"
    "H3"  = "(Vain.) M.Prieto & Wedin"
    "P3"  = "Västra Raikattevuoma Uuijajärvivägen, T lm"
    "Q3"  = 783828.4892504301
    "R3"  = 7532079.29326218
    "AN3" = $null
    "AO3" = $null
    "AW3" = "Mats Williamson"
    "AX3" = "Mats Williamson, * Naturskyddare"

    "A4"  = 475966
    "P4"  = "Raikattivaara V, T lm"
    "Q4"  = 783839.8249707882
    "R4"  = 7531804.227096669
    "AN4" = 1
    "AO4" = "1 substratenheter"
    "AW4" = "Olli Manninen"
    "AX4" = "Olli Manninen, * Naturskyddare"

    "A5"  = 475964
    "B5"  = 89544
    "D5"  = "VU"
    "E5"  = 1503
    "F5"  = "Gräddporing"
    "G5"  = "Sidera lenis"
    "H5"  = "(P.Karst.) Miettinen"
    "Q5"  = 783741.0945520886
    "R5"  = 7532012.121067218

    "A6"  = 675316
    "B6"  = 90840
    "D6"  = "NT"
    "E6"  = 2079
    "F6"  = "Nordtagging"
    "G6"  = "Odonticium romellii"
    "H6"  = "(S.Lundell) Parmasto"
    "Q6"  = 783736.9962432287
    "R6"  = 7532039.150136176

    "A7"  = 475963
    "B7"  = 89544
    "D7"  = "VU"
    "E7"  = 1503
    "F7"  = "Gräddporing"
    "G7"  = "Sidera lenis"
    "H7"  = "(P.Karst.) Miettinen"
    "Q7"  = 783702.3988115358
    "R7"  = 7531995.682737886

    "A8"  = 97482
    "B8"  = 77176
    "D8"  = "NT"
    "E8"  = 353
    "F8"  = "Dvärgbägarlav"
    "G8"  = "Cladonia parasitica"
    "H8"  = "(Hoffm.) Hoffm."
    "Q8"  = 783929.6424552042
    "R8"  = 7532062.489508756

    "A9"  = 475962
    "B9"  = 89544
    "E9"  = 1503
    "F9"  = "Gräddporing"
    "G9"  = "Sidera lenis"
    "H9"  = "(P.Karst.) Miettinen"
    "Q9"  = 783756.0150113704
    "R9"  = 7532029.397363066
}

foreach ($addr in $updates.Keys) {
    $value = $updates[$addr]
    $cell = $ws.Range($addr)
    if ($null -eq $value) {
        $cell.ClearContents()
    } else {
        $cell.Value2 = $value
    }
}
